$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.833.95"
$ws.Range("E2").Value = "  -0.03%  "
$ws.Range("D3").Value = "3.524.42"
$ws.Range("E3").Value = "  +0.50%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "196.56"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.40%  "
$ws.Range("E7").Value = "  -0.55%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.204"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.71%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.647"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.45%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.34"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.40%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000304"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.34%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.49"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.41%  "
$ws.Range("D14").Value = "4.078.50"
$ws.Range("E14").Value = "  +0.24%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "600.27"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.96%  "
$ws.Range("B16").Value = "Uniswap"
$ws.Range("C16").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "12.77"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.42%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "69.960.90"
$ws.Range("E17").Value = "  +0.00%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "19.01"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.50%  "
$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.123"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.96%  "
$ws.Range("B20").Value = "WrappedEther"
$ws.Range("C20").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D20").Value = "3.514.14"
$ws.Range("E20").Value = "  -0.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.984"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.96%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "18.05"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.46%  "
$ws.Range("E23").Value = "  +3.66%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "103.79"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.61"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.22%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.07"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.56%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.80"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.63%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.56"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.81%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.29"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.66%  "
$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.12"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.49%  "
$ws.Range("B31").Value = "dogwifhat"
$ws.Range("C31").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.33"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.37%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.33"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.115"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.16%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.32"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.12%  "
$ws.Range("E35").Value = "  +3.45%  "
$ws.Range("D36").Value = "3.765.74"
$ws.Range("E36").Value = "  +1.47%  "
$ws.Range("D37").Value = "0.0₃0818"
$ws.Range("E37").Value = "  +3.86%  "
$ws.Range("E38").Value = "  +0.04%  "
$ws.Range("E39").Value = "  +0.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.391"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.13%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "505.41"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "36.42"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.87%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.133"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0449"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.41%  "
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.139"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.24%  "
$ws.Range("B46").Value = "ThetaToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.81"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.84%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.24"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.05%  "
$ws.Range("E48").Value = "  +0.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.48"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000251"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.35"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.81%  "
